$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a cell value while preserving it as literal text,
# even when the text looks like a number (e.g. "1.00", "0.0000179").
function Set-TextCell($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "64.033.57"
$ws.Range("E2").Value = "  -1.66%  "
# Row 3
Set-TextCell "D3" "3.503.86"
$ws.Range("E3").Value = "  -0.71%  "
# Row 4
$ws.Range("E4").Value = "  -0.02%  "
# Row 5
Set-TextCell "D5" "584.70"
$ws.Range("E5").Value = "  -1.70%  "
# Row 6
Set-TextCell "D6" "131.49"
$ws.Range("E6").Value = "  -2.39%  "
# Row 7
Set-TextCell "D7" "3.503.81"
$ws.Range("E7").Value = "  -0.58%  "
# Row 8
$ws.Range("E8").Value = "  +0.00%  "
# Row 9
Set-TextCell "D9" "0.483"
$ws.Range("E9").Value = "  -1.90%  "
# Row 10
$ws.Range("E10").Value = "  -0.62%  "
# Row 11
Set-TextCell "D11" "7.11"
$ws.Range("E11").Value = "  -0.20%  "
# Row 12
Set-TextCell "D12" "0.378"
$ws.Range("E12").Value = "  -2.69%  "
# Row 13
Set-TextCell "D13" "4.091.97"
$ws.Range("E13").Value = "  -0.80%  "
# Row 14
Set-TextCell "D14" "27.41"
$ws.Range("E14").Value = "  -1.25%  "
# Row 15
$ws.Range("E15").Value = "  +1.22%  "
# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D16" "0.0000179"
$ws.Range("E16").Value = "  -1.93%  "
# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D17" "3.498.09"
$ws.Range("E17").Value = "  -0.81%  "
# Row 18
Set-TextCell "D18" "64.060.66"
$ws.Range("E18").Value = "  -1.58%  "
# Row 19
Set-TextCell "D19" "9.77"
$ws.Range("E19").Value = "  -3.34%  "
# Row 20
Set-TextCell "D20" "14.11"
$ws.Range("E20").Value = "  -2.50%  "
# Row 21
Set-TextCell "D21" "5.63"
$ws.Range("E21").Value = "  -1.69%  "
# Row 22
Set-TextCell "D22" "382.43"
$ws.Range("E22").Value = "  -2.67%  "
# Row 23
Set-TextCell "D23" "0.571"
$ws.Range("E23").Value = "  -1.55%  "
# Row 24
Set-TextCell "D24" "3.639.42"
$ws.Range("E24").Value = "  -0.87%  "
# Row 25
Set-TextCell "D25" "73.44"
$ws.Range("E25").Value = "  -1.77%  "
# Row 26
$ws.Range("E26").Value = "  -0.07%  "
# Row 27
Set-TextCell "D27" "0.0000116"
$ws.Range("E27").Value = "  +2.41%  "
# Row 28
Set-TextCell "D28" "1.56"
$ws.Range("E28").Value = "  -1.10%  "
# Row 29
Set-TextCell "D29" "7.47"
$ws.Range("E29").Value = "  -3.82%  "
# Row 30
Set-TextCell "D30" "1.00"
$ws.Range("E30").Value = "  +0.06%  "
# Row 31
Set-TextCell "D31" "8.27"
$ws.Range("E31").Value = "  -1.60%  "
# Row 32
$ws.Range("E32").Value = "  -2.20%  "
# Row 33
Set-TextCell "D33" "3.510.01"
$ws.Range("E33").Value = "  -0.64%  "
# Row 34
$ws.Range("E34").Value = "  -0.01%  "
# Row 35
Set-TextCell "D35" "23.41"
$ws.Range("E35").Value = "  -3.35%  "
# Row 36
Set-TextCell "D36" "0.145"
$ws.Range("E36").Value = "  +0.08%  "
# Row 37
Set-TextCell "D37" "5.31"
$ws.Range("E37").Value = "  +0.06%  "
# Row 38
Set-TextCell "D38" "1.56"
$ws.Range("E38").Value = "  -1.68%  "
# Row 39
Set-TextCell "D39" "6.89"
$ws.Range("E39").Value = "  -1.63%  "
# Row 40
Set-TextCell "D40" "158.52"
$ws.Range("E40").Value = "  -6.04%  "
# Row 41
Set-TextCell "D41" "0.0789"
$ws.Range("E41").Value = "  -3.39%  "
# Row 42
Set-TextCell "D42" "0.811"
$ws.Range("E42").Value = "  -1.59%  "
# Row 43
Set-TextCell "D43" "26.28"
$ws.Range("E43").Value = "  +1.58%  "
# Row 44
$ws.Range("E44").Value = "  +0.07%  "
# Row 45
Set-TextCell "D45" "41.94"
$ws.Range("E45").Value = "  -2.35%  "
# Row 46
Set-TextCell "D46" "4.39"
$ws.Range("E46").Value = "  -1.27%  "
# Row 47
$ws.Range("E47").Value = "  -5.75%  "
# Row 48
Set-TextCell "D48" "1.61"
$ws.Range("E48").Value = "  -2.71%  "
# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D49" "2.427.08"
$ws.Range("E49").Value = "  +0.66%  "
# Row 50
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D50" "6.83"
$ws.Range("E50").Value = "  -1.45%  "
# Row 51
Set-TextCell "D51" "0.896"
$ws.Range("E51").Value = "  -1.18%  "
